$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Row 2
$ws1.Range("A2").Value = "Available"
$ws1.Range("B2").Value = "SPA"
$ws1.Range("C2").Value = "BB"
$ws1.Range("D2").Value = 15
$ws1.Range("E2").Value = 1
$ws1.Range("F2").Value = 15
$ws1.Range("G2").Value = 43640.512414733799
$ws1.Range("G2").NumberFormat = "m/d/yy h:mm"

# Row 3
$ws1.Range("A3").Value = "Available"
$ws1.Range("B3").Value = "SPA"
$ws1.Range("C3").Value = "BB"
$ws1.Range("D3").Value = 15
$ws1.Range("E3").Value = 1
$ws1.Range("F3").Value = 1258
$ws1.Range("G3").Value = 43640.512440081016
$ws1.Range("G3").NumberFormat = "m/d/yy h:mm"

# Row 4
$ws1.Range("A4").Value = "Available"
$ws1.Range("B4").Value = "SPA"
$ws1.Range("C4").Value = "BB"
$ws1.Range("D4").Value = 15
$ws1.Range("E4").Value = 1
$ws1.Range("F4").Value = 125
$ws1.Range("G4").Value = 43640.512455370372
$ws1.Range("G4").NumberFormat = "m/d/yy h:mm"

# Row 5
$ws1.Range("A5").Value = "Available"
$ws1.Range("B5").Value = "SPA"
$ws1.Range("C5").Value = "BB"
$ws1.Range("D5").Value = 15
$ws1.Range("E5").Value = 1
$ws1.Range("F5").Value = 158
$ws1.Range("G5").Value = 43640.512481527781
$ws1.Range("G5").NumberFormat = "m/d/yy h:mm"

# Column G width (bestFit) matches "14.85546875" character width from the target file
$ws1.Columns.Item(7).ColumnWidth = 14.85546875

# View: zoom to 115%, selection on K8 (off-data, matches the saved cursor position)
$ws1.Application.ActiveWindow.Zoom = 115
$ws1.Range("K8").Select()
